$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Capture the current (pre-edit) values of rows 411-414, since they will be
# moved down (unchanged) to become rows 415-418 after the new rows are inserted.
$oldRows = @(411, 412, 413, 414)
$savedData = @{}
foreach ($r in $oldRows) {
    $rowVals = @{}
    for ($col = 1; $col -le 18; $col++) {
        $rowVals[$col] = $ws.Cells.Item($r, $col).Value2
    }
    $savedData[$r] = $rowVals
}

# Step 2: Insert 4 new blank rows at 415 so rows 415-476 shift down to 419-480.
$ws.Rows("415:418").Insert()

# Step 3: Write the saved (old) row 411-414 data into the newly inserted rows 415-418.
$destRow = 415
foreach ($r in $oldRows) {
    $rowVals = $savedData[$r]
    for ($col = 1; $col -le 18; $col++) {
        $ws.Cells.Item($destRow, $col).Value = $rowVals[$col]
    }
    $destRow = $destRow + 1
}

# Step 4: Update rows 411-414 with the new week's data.
$ws.Range("D411").Value = 44474
$ws.Range("J411").Value = 120
$ws.Range("K411").Value = 7000
$ws.Range("L411").Value = 8000
$ws.Range("M411").Value = 7500
$ws.Range("P411").Value = 625

$ws.Range("D412").Value = 44474
$ws.Range("J412").Value = 130
$ws.Range("K412").Value = 7000
$ws.Range("L412").Value = 8000
$ws.Range("M412").Value = 7500
$ws.Range("P412").Value = 417

$ws.Range("D413").Value = 44474
$ws.Range("J413").Value = 120
$ws.Range("K413").Value = 9000
$ws.Range("L413").Value = 10000
$ws.Range("M413").Value = 9500
$ws.Range("P413").Value = 792

$ws.Range("D414").Value = 44474
$ws.Range("J414").Value = 130
$ws.Range("K414").Value = 9000
$ws.Range("L414").Value = 10000
$ws.Range("M414").Value = 9500
$ws.Range("P414").Value = 528
